# Updates stock quantities (col F) and recomputed values (col G = Rate * Qty)
# for the affected item rows, plus the dependent "Sub Total" (col B) rows and
# the final grand-total rows, matching the target revision of the
# CryCompanywiseStockReport workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 60
$ws.Range("G9").Value = 1774.2
$ws.Range("B10").Value = 35597.15
$ws.Range("F34").Value = 57
$ws.Range("G34").Value = 1751.04
$ws.Range("B47").Value = 13934.92
$ws.Range("F71").Value = 14
$ws.Range("G71").Value = 3511.06
$ws.Range("F72").Value = 20
$ws.Range("G72").Value = 4183.6
$ws.Range("F73").Value = 6
$ws.Range("G73").Value = 570.12
$ws.Range("F90").Value = 65
$ws.Range("G90").Value = 8771.75
$ws.Range("F92").Value = 49
$ws.Range("G92").Value = 5485.0599999999995
$ws.Range("F97").Value = 329
$ws.Range("G97").Value = 15377.460000000001
$ws.Range("B114").Value = 228595.34
$ws.Range("F195").Value = 112
$ws.Range("G195").Value = 9730.56
$ws.Range("F197").Value = 46
$ws.Range("G197").Value = 2138.54
$ws.Range("F198").Value = 9
$ws.Range("G198").Value = 676.26
$ws.Range("B200").Value = 40853.96
$ws.Range("F250").Value = 9
$ws.Range("G250").Value = 1031.04
$ws.Range("F255").Value = 23
$ws.Range("G255").Value = 737.61
$ws.Range("F260").Value = 13
$ws.Range("G260").Value = 416.91
$ws.Range("F268").Value = 2
$ws.Range("G268").Value = 1110.06
$ws.Range("B274").Value = 62416.83
$ws.Range("B292").Value = 63520
$ws.Range("E292").Value = 153.4
$ws.Range("F292").Value = 69
$ws.Range("G292").Value = 9955.32
$ws.Range("B293").Value = 55373
$ws.Range("E293").Value = 163.62
$ws.Range("F293").Value = -94
$ws.Range("G293").Value = -13562.32
$ws.Range("F310").Value = 7
$ws.Range("G310").Value = 364.91
$ws.Range("B311").Value = 61605
$ws.Range("E311").Value = 133.78
$ws.Range("F311").Value = -13
$ws.Range("G311").Value = -1455.48
$ws.Range("B312").Value = 63563
$ws.Range("E312").Value = 119.04
$ws.Range("F312").Value = 0
$ws.Range("G312").Value = 0
$ws.Range("F321").Value = 21
$ws.Range("G321").Value = 1232.7
$ws.Range("F324").Value = 23
$ws.Range("G324").Value = 1359.99
$ws.Range("F325").Value = 154
$ws.Range("G325").Value = 21224.28
$ws.Range("F328").Value = 139
$ws.Range("G328").Value = 2923.17
$ws.Range("F334").Value = 335
$ws.Range("G334").Value = 50641.95
$ws.Range("F336").Value = 35
$ws.Range("G336").Value = 11218.199999999999
$ws.Range("B339").Value = 247006.28
$ws.Range("F366").Value = 23
$ws.Range("G366").Value = 3158.3599999999997
$ws.Range("F382").Value = 137
$ws.Range("G382").Value = 5892.37
$ws.Range("B395").Value = 222690.41
$ws.Range("F405").Value = 11
$ws.Range("G405").Value = 88.99
$ws.Range("B411").Value = 3024.61
$ws.Range("F426").Value = 29
$ws.Range("G426").Value = 2801.3999999999996
$ws.Range("B430").Value = 37598.35
$ws.Range("F445").Value = 2
$ws.Range("G445").Value = 66.5
$ws.Range("F446").Value = 117
$ws.Range("G446").Value = 7184.969999999999
$ws.Range("B448").Value = 34945.92
$ws.Range("F480").Value = 116
$ws.Range("G480").Value = 1881.52
$ws.Range("F485").Value = 125
$ws.Range("G485").Value = 1643.75
$ws.Range("B487").Value = 64919
$ws.Range("E487").Value = 27.97
$ws.Range("F487").Value = 61
$ws.Range("G487").Value = 1604.3
$ws.Range("B488").Value = 45702
$ws.Range("E488").Value = 31.43
$ws.Range("F488").Value = -215
$ws.Range("G488").Value = -5654.5
$ws.Range("F490").Value = 170
$ws.Range("G490").Value = 2504.1
$ws.Range("B492").Value = -14205.050000000001
$ws.Range("F498").Value = 80
$ws.Range("G498").Value = 4932
$ws.Range("B508").Value = 6418.66
$ws.Range("F548").Value = 181
$ws.Range("G548").Value = 1230.8
$ws.Range("F550").Value = 93
$ws.Range("G550").Value = 1198.77
$ws.Range("F555").Value = 61
$ws.Range("G555").Value = 1991.04
$ws.Range("B556").Value = 4724.97
$ws.Range("B568").Value = 53319
$ws.Range("E568").Value = 310.64
$ws.Range("F568").Value = -6
$ws.Range("G568").Value = -1643.52
$ws.Range("B569").Value = 64810
$ws.Range("E569").Value = 291.22
$ws.Range("F569").Value = 5
$ws.Range("G569").Value = 1369.6000000000001
$ws.Range("F579").Value = 40
$ws.Range("G579").Value = 1060.4
$ws.Range("B586").Value = 18362.23
$ws.Range("B595").Value = 64836
$ws.Range("E595").Value = 104.71
$ws.Range("F595").Value = 0
$ws.Range("G595").Value = 0
$ws.Range("B596").Value = 60031
$ws.Range("E596").Value = 111.69
$ws.Range("F596").Value = -5
$ws.Range("G596").Value = -492.5
$ws.Range("F654").Value = 295
$ws.Range("G654").Value = 23712.1
$ws.Range("B655").Value = 31507.6
$ws.Range("F692").Value = 53
$ws.Range("G692").Value = 13730.71
$ws.Range("B700").Value = 28906.7
$ws.Range("F703").Value = 36
$ws.Range("G703").Value = 1722.96
$ws.Range("F705").Value = 74
$ws.Range("G705").Value = 10591.619999999999
$ws.Range("F710").Value = 18
$ws.Range("G710").Value = 1362.2400000000002
$ws.Range("F711").Value = 9
$ws.Range("G711").Value = 195.48
$ws.Range("F712").Value = 40
$ws.Range("G712").Value = 1491.2
$ws.Range("F714").Value = 57
$ws.Range("G714").Value = 3964.92
$ws.Range("F717").Value = 152
$ws.Range("G717").Value = 20521.519999999997
$ws.Range("F719").Value = 63
$ws.Range("G719").Value = 7604.73
$ws.Range("B720").Value = 55306.939999999995
$ws.Range("F722").Value = 8
$ws.Range("G722").Value = 1995.84
$ws.Range("F724").Value = 7
$ws.Range("G724").Value = 1147.23
$ws.Range("F728").Value = 1
$ws.Range("G728").Value = 150.41
$ws.Range("B736").Value = 65362
$ws.Range("F736").Value = 37
$ws.Range("G736").Value = 1512.1899999999998
$ws.Range("B737").Value = 65079
$ws.Range("F737").Value = 21
$ws.Range("G737").Value = 858.27
$ws.Range("F739").Value = 43
$ws.Range("G739").Value = 1421.5800000000002
$ws.Range("F741").Value = 112
$ws.Range("G741").Value = 4814.88
$ws.Range("F742").Value = 116
$ws.Range("G742").Value = 5787.24
$ws.Range("F745").Value = 69
$ws.Range("G745").Value = 16684.2
$ws.Range("F746").Value = 85
$ws.Range("G746").Value = 4844.150000000001
$ws.Range("B747").Value = 49556.03
$ws.Range("F772").Value = 2608
$ws.Range("G772").Value = 425390.88000000006
$ws.Range("F773").Value = 8
$ws.Range("G773").Value = 1409.36
$ws.Range("F774").Value = 534
$ws.Range("G774").Value = 151052.58000000002
$ws.Range("F775").Value = 417
$ws.Range("G775").Value = 60319.05
$ws.Range("F778").Value = 50
$ws.Range("G778").Value = 6429.000000000001
$ws.Range("B779").Value = 645235.54
$ws.Range("F782").Value = 34
$ws.Range("G782").Value = 4964.34
$ws.Range("B796").Value = 60947.979999999996
$ws.Range("B797").Value = 2403721.29
$ws.Range("B798").Value = 2403721.29
